$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.674.08"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.916.17"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.52"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.38"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.15"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0885"
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.69"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.89"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.380.82"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.922.09"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.976"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.712.45"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.53"
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.26"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.91"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0980"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.77"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.11"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("E26").Value = "  +12.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.04"
$ws.Range("E27").Value = "  +2.67%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  +14.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.106"
$ws.Range("E30").Value = "  +13.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.55"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.40"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.03"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "52.28"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0439"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -16.15%  "
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.34"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.97"
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.92"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("E47").Value = "  -4.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.124.41"
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.248"
$ws.Range("E49").Value = "  -7.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0332"
$ws.Range("E50").Value = "  +4.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.08"
$ws.Range("E51").Value = "  -0.71%  "
